$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'241.28"
$ws.Range("D3").Value = "'21.43"
$ws.Range("D4").Value = "'5.171"
$ws.Range("D5").Value = "'0.05529"
$ws.Range("D6").Value = "'3.370"
$ws.Range("D7").Value = "'6.327"
$ws.Range("D8").Value = "'0.8039"
$ws.Range("D9").Value = "'0.9534"
$ws.Range("D10").Value = "'0.1378"
$ws.Range("D11").Value = "'0.07309"
$ws.Range("D12").Value = "'0.03021"
$ws.Range("D13").Value = "'0.03068"
$ws.Range("D14").Value = "'0.09299"
$ws.Range("D15").Value = "'3.585"
$ws.Range("D16").Value = "'0.001655"
$ws.Range("D17").Value = "'0.04680"
$ws.Range("D18").Value = "'0.0005766"
$ws.Range("D19").Value = "'0.006402"
$ws.Range("D20").Value = "'0.004980"
$ws.Range("D21").Value = "'0.001044"
$ws.Range("D24").Value = "'3.771"
$ws.Range("D26").Value = "'0.3239"
$ws.Range("D40").Value = "'0.03824"
$ws.Range("D41").Value = "'0.006910"
$ws.Range("D42").Value = "'0.1027"
$ws.Range("D43").Value = "'0.003097"
$ws.Range("D44").Value = "'0.008279"
$ws.Range("D45").Value = "'0.00005962"
$ws.Range("D47").Value = "'0.0005515"
$ws.Range("D49").Value = "'0.1096"
